$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("A2").Value = 44207
$ws.Range("B2").Value = 44216
$ws.Range("C2").Value = 22534.30314723
$ws.Range("D2").Value = 18227.88293334
$ws.Range("E2").Value = 5015.074000000001
$ws.Range("F2").Value = 5015.074000000001
$ws.Range("G2").Value = 19047.28366525
$ws.Range("H2").Value = 25657.87662308999
$ws.Range("I2").Value = 0.1547427253107065
